$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 96

# A96: date value with same style as the date column (A95) -> yyyy-mm-dd hh:mm:ss
$cA = $ws.Cells.Item($row, 1)
$cA.Font.Name = "Calibri"
$cA.Font.Size = 11
$cA.NumberFormat = "yyyy-mm-dd hh:mm:ss"
$cA.Value = 45453.2916666667

# B96-F96: plain numbers
$ws.Cells.Item($row, 2).Value = 0
$ws.Cells.Item($row, 3).Value = 2
$ws.Cells.Item($row, 4).Value = 2
$ws.Cells.Item($row, 5).Value = 2
$ws.Cells.Item($row, 6).Value = 2

# G96: text "2" (shared string), reset formatting back to default afterwards
$cG = $ws.Cells.Item($row, 7)
$cG.NumberFormat = "@"
$cG.Value = "2"
$cG.Style = "Normal"

# H96: ticker text "KK.MI"
$ws.Cells.Item($row, 8).Value = "KK.MI"
